$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Replace the text category headers (B1:E1 = "cardinality: N") on every
#    sheet with plain numbers (10 / 100 / 1000 / 10000). This turns the
#    cells from shared-string references into numeric cells, which also
#    drops the four now-unused "cardinality: N" entries from the shared
#    string table automatically.
# ---------------------------------------------------------------------------
$sheetNames = @("uniform", "normal", "exponential 1", "exponential 2")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = 10
    $ws.Range("C1").Value = 100
    $ws.Range("D1").Value = 1000
    $ws.Range("E1").Value = 10000
}

# ---------------------------------------------------------------------------
# 2) Update each chart's title so the sheet the chart belongs to is called
#    out explicitly instead of the generic "per cardinality" wording.
# ---------------------------------------------------------------------------
$wsUniform = $wb.Worksheets.Item("uniform")
$chartUniform = $wsUniform.ChartObjects().Item(1).Chart
$chartUniform.ChartTitle.Text = "Encoding performance [uniform]"

$wsNormal = $wb.Worksheets.Item("normal")
$chartNormal = $wsNormal.ChartObjects().Item(1).Chart
$chartNormal.ChartTitle.Text = "Encoding performance [uniform]"

$wsExp1 = $wb.Worksheets.Item("exponential 1")
$chartExp1 = $wsExp1.ChartObjects().Item(1).Chart
$chartExp1.ChartTitle.Text = "Encoding performance [exponential 1]"

$wsExp2 = $wb.Worksheets.Item("exponential 2")
$chartExp2 = $wsExp2.ChartObjects().Item(1).Chart
$chartExp2.ChartTitle.Text = "Encoding performance [exponential 2]"

# ---------------------------------------------------------------------------
# 3) Update the selected cell on every sheet and make "exponential 2" the
#    active (visible) tab, mirroring the saved UI state in the workbook.
# ---------------------------------------------------------------------------
$wsUniform.Range("G29").Select() | Out-Null
$wsNormal.Range("G29").Select() | Out-Null
$wsExp1.Range("I24").Select() | Out-Null
$wsExp2.Range("G30").Select() | Out-Null
$wsExp2.Activate() | Out-Null
